$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The quarterly report rolled forward by one quarter: the oldest quarter
# (column D, "1399-09-30 (2)") drops off, every remaining quarter shifts one
# column to the left, and a brand new quarter ("1402-02-28 (8)") is appended
# as the new last column (M).
# ---------------------------------------------------------------------------

# Drop the oldest quarter column - this shifts D:M left to C:L automatically.
$ws.Columns("D").Delete()

# Bring column M to life with the same formatting as column L (the previous
# last data column), then give it its own (slightly wider) column width.
$ws.Range("L1:L59").Copy()
$ws.Range("M1:M59").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Columns("M").ColumnWidth = 30.1

# Company name correction.
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# New quarter header label (row 8) and publish-date row (row 9).
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("I9").Value = "1402-02-28 (8)"
# "1402-02-28" alone parses as a (Persian-calendar) date literal, which Excel
# would silently convert to a numeric serial; force it to stay text, then
# reapply row 9's formatting (the forced-text entry nudges the style).
$ws.Range("M9").Value = "'1402-02-28"
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New quarter's balance-sheet figures (column M).
$ws.Range("M12").Value = 212241
$ws.Range("M13").Value = 238637
$ws.Range("M14").Value = 5429174
$ws.Range("M15").Value = 913837
$ws.Range("M16").Value = 101223
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 6895112
$ws.Range("M19").Value = 0
$ws.Range("M20").Value = 1461478
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 361220
$ws.Range("M23").Value = 51030
$ws.Range("M24").Value = 40
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 1873728
$ws.Range("M27").Value = 8768840
$ws.Range("M29").Value = 1522103
$ws.Range("M30").Value = 40
$ws.Range("M31").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("M33").Value = 142475
$ws.Range("M34").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M37").Value = 1664578
$ws.Range("M38").Value = 0
$ws.Range("M39").Value = 40
$ws.Range("M40").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("M43").Value = 1664578
$ws.Range("M45").Value = 3000000
$ws.Range("M46").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("M48").Value = -235820
$ws.Range("M49").Value = 115430
$ws.Range("M50").Value = 294770
$ws.Range("M51").Value = 0
$ws.Range("M52").Value = 40
$ws.Range("M53").Value = 0
$ws.Range("M54").Value = 40
$ws.Range("M55").Value = 0
$ws.Range("M56").Value = 3929882
$ws.Range("M57").Value = 7104262
$ws.Range("M58").Value = 8768840
